$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("M4").Value = 2093.72
$ws1.Range("M8").Value = 663.55
$ws1.Range("D10").Value = 4138.56
$ws1.Range("L10").Value = -297.88
$ws1.Range("L17").Value = 3162.93

# Totals row (counts of "X de 24")
$ws1.Range("D26").Value = "1 de 24"
$ws1.Range("L26").Value = "2 de 24"
$ws1.Range("M26").Value = "6 de 24"

# --- Sheet 2: "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F4").Value = 2093.72
$ws2.Range("F8").Value = 663.55
$ws2.Range("F10").Value = 3840.68
$ws2.Range("F17").Value = 3162.93
$ws2.Range("F26").Value = 36711.93

# --- Sheet 3: "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D3").Value = 4138.56
$ws3.Range("E3").Value = 3336.8483879616
$ws3.Range("F3").Value = 0.5536232651402349

$ws3.Range("D11").Value = 4477.07
$ws3.Range("E11").Value = -1554.84541814726
$ws3.Range("F11").Value = 1.532075949193974

$ws3.Range("D12").Value = 25372.73
$ws3.Range("E12").Value = 2582.25
$ws3.Range("F12").Value = 0.9076282651606261

$ws3.Range("D14").Value = 35519.04
$ws3.Range("E14").Value = 6684.341100094684
$ws3.Range("F14").Value = 0.841615981329996
